# total_hits.xlsx — append the next slate of games (2022-05-03, serial 44684)
# to the bottom of Sheet1, matching the existing table's layout/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently runs through row 261. New rows start at 262.
$startRow = 262

# Give column A (the date column) the same number-formatted style ("s") as
# the row above it, so new date cells stay formatted like the rest of the
# column instead of getting a brand-new style entry.
$ws.Range("A261").Copy() | Out-Null
$ws.Range("A262:A277").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Each entry: date serial, visitor (B), home (C), and optional line/over/under
# (D/E/F) — games that haven't posted a total yet just omit D/E/F, same as
# the existing rows for in-progress/unlisted games.
$rows = @(
    @{ B = "Atlanta Braves";       C = "New York Mets" },
    @{ B = "Atlanta Braves";       C = "New York Mets" },
    @{ B = "Arizona Diamondbacks"; C = "Miami Marlins" },
    @{ B = "Cinncinatti Reds";     C = "Milwaukee Brewers";    D = 15.5; E = 105;  F = -135 },
    @{ B = "Washington Nationals"; C = "Colorado Rockies";     D = 18.5; E = -105; F = -135 },
    @{ B = "San Francisco Giants"; C = "Los Angeles Dodgers";  D = 14.5; E = -140; F = 105 },
    @{ B = "Minnesota Twins";      C = "Baltimore Orioles";    D = 15.5; E = 105;  F = -140 },
    @{ B = "New York Yankees";     C = "Toronto Blue Jays";    D = 16.5; E = 125;  F = -165 },
    @{ B = "Los Angeles Angels";   C = "Boston Red Sox";       D = 16.5; E = -110; F = -125 },
    @{ B = "Seattle Mariners";     C = "Houston Astros";       D = 16.5; E = -110; F = -125 },
    @{ B = "Tampa Bay Rays";       C = "Oakland Athletics";    D = 15.5; E = 105;  F = -135 },
    @{ B = "San Diego Padres";     C = "Cleveland Gaurdians" },
    @{ B = "Texas Rangers";        C = "Philidelphia Phillies"; D = 16.5; E = 105;  F = -145 },
    @{ B = "Pittsburgh Pirates";   C = "Detroit Tigers";       D = 16.5; E = -110; F = -125 },
    @{ B = "Chicago White Sox";    C = "Chicago Cubs";         D = 15.5; E = 120;  F = -160 },
    @{ B = "Saint Luis Cardinals"; C = "Kansas City Royals";   D = 15.5; E = 105;  F = -135 }
)

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = 44684
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    if ($row.ContainsKey("D")) {
        $ws.Cells.Item($r, 4).Value = $row.D
        $ws.Cells.Item($r, 5).Value = $row.E
        $ws.Cells.Item($r, 6).Value = $row.F
    }
    $r = $r + 1
}

# Leave the selection where a user would land after typing the last row.
$lastRow = $r
$ws.Cells.Item($lastRow, 6).Select() | Out-Null
